# Add experiment 7 and 8 results: fill in the min/mean/max message error
# percent columns (G, H, I) for sheet exp_7, rows 2-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("exp_7")

$data = @(
    @(2,  0,                    0,                    0),
    @(3,  0,                    0,                    0),
    @(4,  0,                    0,                    0),
    @(5,  0,                    1.1764705882352899,   1.9607843137254899),
    @(6,  0,                    1.2,                  4),
    @(7,  0,                    1.19215686274509,     4),
    @(8,  0,                    3.73700088731144,     6.5217391304347698),
    @(9,  0,                    0,                    0),
    @(10, 1.9607843137254899,   1.9921568627451001,   2),
    @(11, 0,                    0.81632653061224503,  4.0816326530612201),
    @(12, 0,                    1.60032012805122,     4),
    @(13, 0,                    0,                    0),
    @(14, 0,                    0.39215686274509798,  1.9607843137254899),
    @(15, 0,                    0.39215686274509798,  1.9607843137254899),
    @(16, 0,                    1.97647058823529,     4),
    @(17, 0,                    1.9921568627451001,   4),
    @(18, 0,                    1.19215686274509,     4),
    @(19, 0,                    5.8576752440106397,   13.043478260869501),
    @(20, 0,                    1.5610859728506801,   1.9607843137254899),
    @(21, 1.9607843137254899,   3.5921568627450999,   4),
    @(22, 1.9607843137254899,   3.6088035214085599,   6.1224489795918302),
    @(23, 0,                    2.80880352140856,     6.1224489795918302)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 7).Value = $row[1]
    $ws.Cells.Item($r, 8).Value = $row[2]
    $ws.Cells.Item($r, 9).Value = $row[3]
}

# Make exp_7 the active sheet/tab and set its selection, matching the
# workbook's activeTab and the sheet's tabSelected/selection state.
$ws.Activate()
$ws.Range("K10").Select()
